# Added remote webdriver support
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("DATAMANAGER")

# Insert a new column E ("mode") into DATAMANAGER, before the existing
# browser/url columns, and fill in the mode values.
$wsData.Columns.Item(5).Insert()

$wsData.Range("E1").Value = "mode"
$wsData.Range("E2").Value = "remote"
$wsData.Range("E3").Value = "local"
$wsData.Range("E4").Value = "local"
$wsData.Range("E5").Value = "local"

# Make DATAMANAGER the active/selected sheet, with E2 as the selected cell.
$wsData.Activate()
$wsData.Range("E2").Select()
